# Fruta / hortaliza, semanal
# Weekly refresh: re-shuffle the Fecha/Volumen/Precio columns (D, J, K, L, M, P)
# across the data rows (2-27) of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, in column order: D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$rowData = @{
    2  = @(44755, 50, 20000, 20000, 20000, 1333)
    3  = @(44830, 25, 12000, 12000, 12000, 800)
    4  = @(44826, 50, 20000, 20000, 20000, 1333)
    5  = @(44525, 40, 8000, 8000, 8000, 533)
    6  = @(44767, 50, 20000, 20000, 20000, 1333)
    7  = @(44813, 20, 20000, 20000, 20000, 1333)
    8  = @(44769, 50, 20000, 20000, 20000, 1333)
    9  = @(44837, 80, 16000, 16000, 16000, 1067)
    10 = @(44749, 50, 20000, 20000, 20000, 1333)
    11 = @(44811, 30, 20000, 20000, 20000, 1333)
    12 = @(44518, 50, 10000, 10000, 10000, 667)
    13 = @(44756, 80, 20000, 20000, 20000, 1333)
    14 = @(44845, 20, 16000, 16000, 16000, 1067)
    15 = @(44825, 30, 20000, 20000, 20000, 1333)
    16 = @(44812, 80, 20000, 20000, 20000, 1333)
    17 = @(44757, 30, 20000, 20000, 20000, 1333)
    18 = @(44841, 20, 16000, 16000, 16000, 1067)
    19 = @(44838, 10, 20000, 20000, 20000, 1333)
    20 = @(44776, 80, 20000, 20000, 20000, 1333)
    21 = @(44771, 40, 20000, 20000, 20000, 1333)
    22 = @(44839, 80, 16000, 16000, 16000, 1067)
    23 = @(45134, 5, 20000, 20000, 20000, 1333)
    24 = @(44827, 20, 20000, 20000, 20000, 1333)
    25 = @(44824, 20, 20000, 20000, 20000, 1333)
    26 = @(44819, 100, 20000, 20000, 20000, 1333)
    27 = @(44508, 40, 10000, 10000, 10000, 667)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
    $ws.Cells.Item($r, 11).Value = $vals[2]
    $ws.Cells.Item($r, 12).Value = $vals[3]
    $ws.Cells.Item($r, 13).Value = $vals[4]
    $ws.Cells.Item($r, 16).Value = $vals[5]
}
